# edit.ps1
# 1. Update the cached "datetimeFigureOut" footer-date text from 9/11/22 to
#    12/14/22 on the slide master and every slide layout.
# 2. Remove the old "TextBox 45" shape (hyperlink to the hyperledgendary repo)
#    from slide 1.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "9/11/22") {
                $shp.TextFrame.TextRange.Text = "12/14/22"
            }
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every slide layout (custom layout) attached to the master ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- Remove the obsolete hyperlink textbox from slide 1 ---
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 45") {
        $shp.Delete()
    }
}
